$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '70.496.94'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +1.63%  '
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.561.60'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +0.68%  '
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '611.42'
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '173.16'
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +0.75%  '
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.618'
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +1.10%  '
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '3.558.34'
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c.ClearFormats()

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  +2.81%  '
$c.ClearFormats()

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '7.63'
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +14.04%  '
$c.ClearFormats()

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -2.06%  '
$c.ClearFormats()

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +0.57%  '
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.143.57'
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +1.11%  '
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '8.39'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -2.07%  '
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '612.37'
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -2.45%  '
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.559.97'
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +0.59%  '
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '70.591.17'
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +1.78%  '
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -2.25%  '
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.41'
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -0.75%  '
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.885'
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.76%  '
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.43'
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -15.98%  '
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '16.08'
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '97.02'
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -0.87%  '
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -0.79%  '
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '33.46'
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +1.61%  '
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '9.08'
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -2.55%  '
$c.ClearFormats()

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.73%  '
$c.ClearFormats()

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -3.37%  '
$c.ClearFormats()

$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = 'NEARProtocol'
$c.ClearFormats()
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '6.99'
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -0.29%  '
$c.ClearFormats()

$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = 'Mantle'
$c.ClearFormats()
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.30'
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -2.55%  '
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '576.68'
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -9.54%  '
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.66'
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +4.29%  '
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -1.92%  '
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +0.19%  '
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0483'
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +6.10%  '
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '57.40'
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.ClearFormats()

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +0.06%  '
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.385.67'
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -0.28%  '
$c.ClearFormats()

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -3.00%  '
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '33.24'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +0.65%  '
$c.ClearFormats()

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +7.35%  '
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +1.34%  '
$c.ClearFormats()

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '133.93'
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +1.30%  '
$c.ClearFormats()
